$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.185.12'
$ws.Range('D2').Style = $origStyle
$ws.Range('E2').Value = '  +0.16%  '

$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.903.08'
$ws.Range('D3').Style = $origStyle

$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = $origStyle
$ws.Range('E4').Value = '  -0.05%  '

$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.31'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -0.69%  '

$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +0.00%  '

$ws.Range('E7').Value = '  +1.36%  '

$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3774'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  +1.36%  '

$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07259'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  +0.53%  '

$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.12'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +0.40%  '

$ws.Range('E11').Value = '  -0.75%  '

$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08391'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +9.89%  '

$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.900.83'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +0.30%  '

$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.62'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -0.38%  '

$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.266'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  -0.39%  '

$ws.Range('E16').Value = '  -0.08%  '

$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008611'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +1.12%  '

$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.56'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +1.27%  '

$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +0.00%  '

$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.225.17'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +0.20%  '

$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.057'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -0.12%  '

$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.161.24'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +1.29%  '

$ws.Range('E23').Value = '  -0.37%  '

$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.441'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  -0.11%  '

$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.81'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  +0.90%  '

$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.278'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +5.63%  '

$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.755'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -2.19%  '

$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.17'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +0.40%  '

$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.73'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -0.06%  '

$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.932'
$ws.Range('D30').Style = $origStyle

$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.799'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  -0.39%  '

$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09280'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +0.71%  '

$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8093'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +6.12%  '

$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05066'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  -0.05%  '

$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.235'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +3.55%  '

$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.951'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  -2.66%  '

$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.352'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +2.05%  '

$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.613'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +2.23%  '

$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5696'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +1.21%  '

$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01990'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  -0.36%  '

$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.070'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -0.80%  '

$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.652'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +0.70%  '

$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.964'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +0.41%  '

$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '117.97'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -0.85%  '

$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1513'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +0.11%  '

$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4844'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +0.66%  '

$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.000'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -0.02%  '

$ws.Range('E48').Value = '  -0.20%  '

$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.613'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +2.35%  '

$ws.Range('E50').Value = '  +0.61%  '

$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.64'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -0.02%  '
